$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

$ws.Range("C6").Value = "TestSBO_Replay_M1/Activity/ActivityData.xlsx"
$ws.Range("C9").Value = "TestSBO_Replay_M1/Activity/ActivityData_D2.xlsx"
$ws.Range("C12").Value = "TestSBO_Replay_M1/Activity/ActivityData_D3.xlsx"
$ws.Range("C14").Value = "TestSBO_Replay_M1/Activity/ActivityData_D4.xlsx"
$ws.Range("C18").Value = "TestSBO_Replay_M1/Activity/ActivityData_D5.xlsx"

$ws.Activate()
$ws.Range("C13").Select()
